# ==========================================================================
# Applies the "Add files via upload" commit to the workbook:
#   - NutritionalData sheet (sheet1): row 237 formulas re-pointed from
#     B331 -> B329; row 247 rebuilt with a new ingredient (haribo gummy
#     bears) and its own formulas; six brand-new ingredient rows 248-253
#     appended (hersheys bar, rose wine, pita bread, ponzu sauce, double
#     salmon poki bowl, and a creamcheese-doubling row).
#   - researchMeasures sheet (sheet2): day row 95 gets its ingredient
#     list + nutrition totals filled in (falafel-platter day) and a new
#     day row 96 is appended (poki bowl day).
# ==========================================================================

$wb = $excel.ActiveWorkbook
$wsNutrition = $wb.Worksheets.Item("NutritionalData")
$wsResearch  = $wb.Worksheets.Item("researchMeasures")

# --------------------------------------------------------------------
# NutritionalData!A237:H237 - re-point the shared formula from B331 to
# B329 (now written out per-cell, matching the no-longer-shared result)
# --------------------------------------------------------------------
foreach ($col in @("B","C","D","E","F","G","H")) {
    $wsNutrition.Range($col + "237").Formula = "=SUM(" + $col + "40*2," + $col + "39," + $col + "37*4.5," + $col + "329*5," + $col + "115*5)"
}

# --------------------------------------------------------------------
# NutritionalData!A247:H247 - replace with "haribo gummy bears"
# --------------------------------------------------------------------
$wsNutrition.Range("A247").Value = "haribo gummy bears 9 servings"
$wsNutrition.Range("A230").Copy()
$wsNutrition.Range("A247").PasteSpecial(-4122)

$wsNutrition.Range("B247").Formula = "=100*9"
$wsNutrition.Range("C247").Formula = "=0*9"
$wsNutrition.Range("D247").Formula = "=0*9"
$wsNutrition.Range("E247").Formula = "=2*9"
$wsNutrition.Range("F247").Formula = "=23*9"
$wsNutrition.Range("G247").Formula = "=14*9"
$wsNutrition.Range("H247").Formula = "=5*9"

# --------------------------------------------------------------------
# NutritionalData!A248:H248 - "hersheys chocolate almond and toffee bar"
# --------------------------------------------------------------------
$wsNutrition.Range("A248").Value = "hersheys chocolate almond and toffee bar 4 servings per Xl bar"
$wsNutrition.Range("A230").Copy()
$wsNutrition.Range("A248").PasteSpecial(-4122)

$wsNutrition.Range("B248").Formula = "=150*4"
$wsNutrition.Range("C248").Formula = "=9*4"
$wsNutrition.Range("D248").Formula = "=5*4"
$wsNutrition.Range("E248").Formula = "=3*4"
$wsNutrition.Range("F248").Formula = "=17*4"
$wsNutrition.Range("G248").Formula = "=1*4"
$wsNutrition.Range("H248").Formula = "=50*4"

# --------------------------------------------------------------------
# NutritionalData!A249:H249 - rose wine bottle
# --------------------------------------------------------------------
$wsNutrition.Range("A249").Value = "1 BOTTLE WINE ROSE GENERIC 3 GLASSES, https://www.calorieking.com/us/en/foods/f/calories-in-wines-cabernet-franc-red-wine-10-alc/-kpHuBQISOKXQPyfyKIobA"
$wsNutrition.Range("A230").Copy()
$wsNutrition.Range("A249").PasteSpecial(-4122)

$wsNutrition.Range("B241").Copy()
$wsNutrition.Range("B249").PasteSpecial(-4122)
$wsNutrition.Range("B249").Value = 377
$wsNutrition.Range("C249").Value = 0
$wsNutrition.Range("D249").Value = 0
$wsNutrition.Range("E249").Value = 0.3
$wsNutrition.Range("F249").Value = 10.9
$wsNutrition.Range("G249").Formula = "=G49*2"
$wsNutrition.Range("H249").Value = 0
$wsNutrition.Range("C249:H249").NumberFormat = "0.00"

# --------------------------------------------------------------------
# NutritionalData!A250:H250 - pita bread slice
# --------------------------------------------------------------------
$wsNutrition.Range("A250").Value = "1 slice pita bread, https://www.calorieking.com/us/en/foods/f/calories-in-bread-rolls-buns-white-pita-khoubiz-pocket/gVZwM5UNRSGaUF9qX5K2qg"
$wsNutrition.Range("A230").Copy()
$wsNutrition.Range("A250").PasteSpecial(-4122)

$wsNutrition.Range("B250").Value = 165
$wsNutrition.Range("C250").Value = 0.7
$wsNutrition.Range("D250").Value = 0.1
$wsNutrition.Range("E250").Value = 5.5
$wsNutrition.Range("F250").Value = 33.4
$wsNutrition.Range("G250").Value = 1.3
$wsNutrition.Range("H250").Value = 322

# --------------------------------------------------------------------
# NutritionalData!A251:H251 - ponzu poki bowl sauce
# --------------------------------------------------------------------
$wsNutrition.Range("A251").Value = "ponzu poki bowl sauce, https://www.calorieking.com/us/en/foods/f/calories-in-sauces-ponzu-or-ponzu-lime/7IfSf4aVS4WoaU11NMCwnA"
$wsNutrition.Range("A230").Copy()
$wsNutrition.Range("A251").PasteSpecial(-4122)

$wsNutrition.Range("B251").Value = 40
$wsNutrition.Range("C251").Value = 0
$wsNutrition.Range("D251").Value = 0
$wsNutrition.Range("E251").Value = 1
$wsNutrition.Range("F251").Value = 8
$wsNutrition.Range("G251").Value = 0
$wsNutrition.Range("H251").Value = 1520

# --------------------------------------------------------------------
# NutritionalData!A252:H252 - double salmon poki bowl
# --------------------------------------------------------------------
$wsNutrition.Range("A252").Value = "double salmon poki bowl with ponzu sauce instead of teriyaki sauce in cell A177 ingredients"
$wsNutrition.Range("A230").Copy()
$wsNutrition.Range("A252").PasteSpecial(-4122)

$wsNutrition.Range("B252").Value = 578.5
$wsNutrition.Range("C252").Value = 11.074999999999999
$wsNutrition.Range("D252").Value = 2.15
$wsNutrition.Range("E252").Value = 17.699999999999996
$wsNutrition.Range("F252").Value = 101.75
$wsNutrition.Range("G252").Value = 8.25
$wsNutrition.Range("H252").Value = 1731

# --------------------------------------------------------------------
# NutritionalData!B253:H253 - creamcheese 1/4 cup (double row 193)
# --------------------------------------------------------------------
$wsNutrition.Range("B241").Copy()
$wsNutrition.Range("B253:H253").PasteSpecial(-4122)
foreach ($col in @("B","C","D","E","F","G","H")) {
    $wsNutrition.Range($col + "253").Formula = "=" + $col + "193*2"
}

Write-Output "NutritionalData rows done"
